# Auto-generated edits applying the Ultima_Profits market-price refresh diff.
# Updates columns H-N (price/profit figures) for specific leve rows across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4873.273
$ws.Range("I76").Value = 2950.75
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 2950.75
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -2635.75
$ws.Range("N76").Value = -10630

$ws.Range("H79").Value = 4873.273
$ws.Range("I79").Value = 2950.75
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 2950.75
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -1858.75
$ws.Range("N79").Value = -12184

$ws.Range("H118").Value = 926.61536
$ws.Range("I118").Value = 318.5
$ws.Range("J118").Value = 1899.6
$ws.Range("K118").Value = 955.5
$ws.Range("L118").Value = 5698.799999999999
$ws.Range("M118").Value = 701.5
$ws.Range("N118").Value = -9012.799999999999

$ws.Range("H129").Value = 1869.1526
$ws.Range("I129").Value = 398
$ws.Range("K129").Value = 1194
$ws.Range("M129").Value = 3806

$ws.Range("H132").Value = 5686417.5
$ws.Range("I132").Value = 4520.4326
$ws.Range("J132").Value = 35719304
$ws.Range("K132").Value = 13561.2978
$ws.Range("L132").Value = 107157912
$ws.Range("M132").Value = -11031.2978
$ws.Range("N132").Value = -107162972

$ws.Range("H137").Value = 11112179
$ws.Range("I137").Value = 682.2
$ws.Range("K137").Value = 2046.6
$ws.Range("M137").Value = 503.3999999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8570.18
$ws.Range("I32").Value = 8451.407999999999
$ws.Range("J32").Value = 9038.666999999999
$ws.Range("K32").Value = 8451.407999999999
$ws.Range("L32").Value = 9038.666999999999
$ws.Range("M32").Value = -8164.407999999999
$ws.Range("N32").Value = -9612.666999999999

$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976

$ws.Range("H45").Value = 1622
$ws.Range("I45").Value = 1552.6666
$ws.Range("K45").Value = 1552.6666
$ws.Range("M45").Value = -1175.6666

$ws.Range("H74").Value = 6850989
$ws.Range("I74").Value = 8621877
$ws.Range("J74").Value = 3555.2
$ws.Range("K74").Value = 8621877
$ws.Range("L74").Value = 3555.2
$ws.Range("M74").Value = -8621003
$ws.Range("N74").Value = -5303.2

$ws.Range("H77").Value = 6850989
$ws.Range("I77").Value = 8621877
$ws.Range("J77").Value = 3555.2
$ws.Range("K77").Value = 43109385
$ws.Range("L77").Value = 17776
$ws.Range("M77").Value = -43105017
$ws.Range("N77").Value = -26512

$ws.Range("H132").Value = 5210201.5
$ws.Range("I132").Value = 6758476
$ws.Range("J132").Value = 2368.7273
$ws.Range("K132").Value = 20275428
$ws.Range("L132").Value = 7106.1819
$ws.Range("M132").Value = -20272898
$ws.Range("N132").Value = -12166.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 19000
$ws.Range("J39").Value = 19000
$ws.Range("L39").Value = 19000
$ws.Range("N39").Value = -19778

$ws.Range("H134").Value = 2511.3655
$ws.Range("I134").Value = 1358.1714
$ws.Range("J134").Value = 4885.5884
$ws.Range("K134").Value = 4074.5142
$ws.Range("L134").Value = 14656.7652
$ws.Range("M134").Value = -1539.5142
$ws.Range("N134").Value = -19726.7652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5559512.5
$ws.Range("I31").Value = 4002.4443
$ws.Range("J31").Value = 22226042
$ws.Range("K31").Value = 4002.4443
$ws.Range("L31").Value = 22226042
$ws.Range("M31").Value = -3707.4443
$ws.Range("N31").Value = -22226632

$ws.Range("H34").Value = 5559512.5
$ws.Range("I34").Value = 4002.4443
$ws.Range("J34").Value = 22226042
$ws.Range("K34").Value = 4002.4443
$ws.Range("L34").Value = 22226042
$ws.Range("M34").Value = -3800.4443
$ws.Range("N34").Value = -22226446

$ws.Range("H54").Value = 30046
$ws.Range("J54").Value = 30046
$ws.Range("L54").Value = 30046
$ws.Range("N54").Value = -31362

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null

$ws.Range("H122").Value = 3670.6667
$ws.Range("I122").Value = 3670.6667
$ws.Range("K122").Value = 11012.0001
$ws.Range("M122").Value = -8562.000100000001

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null

$ws.Range("H132").Value = 7354258.5
$ws.Range("I132").Value = 8773154
$ws.Range("J132").Value = 1797.4546
$ws.Range("K132").Value = 26319462
$ws.Range("L132").Value = 5392.3638
$ws.Range("M132").Value = -26316932
$ws.Range("N132").Value = -10452.3638

$ws.Range("H134").Value = 582188.25
$ws.Range("I134").Value = 1523.1818
$ws.Range("J134").Value = 2977431.5
$ws.Range("K134").Value = 4569.5454
$ws.Range("L134").Value = 8932294.5
$ws.Range("M134").Value = -2034.5454
$ws.Range("N134").Value = -8937364.5

$ws.Range("H140").Value = 39516.668
$ws.Range("J140").Value = 39516.668
$ws.Range("L140").Value = 39516.668
$ws.Range("N140").Value = -49876.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 631374
$ws.Range("I2").Value = 89.61539
$ws.Range("J2").Value = 1543229.2
$ws.Range("K2").Value = 537.6923400000001
$ws.Range("L2").Value = 9259375.199999999
$ws.Range("M2").Value = -424.6923400000001
$ws.Range("N2").Value = -9259601.199999999

$ws.Range("H3").Value = 5885.636
$ws.Range("I3").Value = 3907
$ws.Range("J3").Value = 6909.069
$ws.Range("K3").Value = 11721
$ws.Range("L3").Value = 20727.207
$ws.Range("M3").Value = -11609
$ws.Range("N3").Value = -20951.207

$ws.Range("H34").Value = 1213.7222
$ws.Range("I34").Value = 168.09091
$ws.Range("J34").Value = 2856.8572
$ws.Range("K34").Value = 504.27273
$ws.Range("L34").Value = 8570.571599999999
$ws.Range("M34").Value = -420.27273
$ws.Range("N34").Value = -8738.571599999999

$ws.Range("H42").Value = 4332.3335
$ws.Range("I42").Value = 1000
$ws.Range("J42").Value = 4998.8
$ws.Range("K42").Value = 3000
$ws.Range("L42").Value = 14996.4
$ws.Range("M42").Value = -2466
$ws.Range("N42").Value = -16064.4

$ws.Range("H121").Value = 1141.5
$ws.Range("I121").Value = 250
$ws.Range("J121").Value = 1240.5555
$ws.Range("K121").Value = 750
$ws.Range("L121").Value = 3721.6665
$ws.Range("M121").Value = 560
$ws.Range("N121").Value = -6341.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7874
$ws.Range("I70").Value = 10877.071
$ws.Range("J70").Value = 4370.4165
$ws.Range("K70").Value = 10877.071
$ws.Range("L70").Value = 4370.4165
$ws.Range("M70").Value = -10607.071
$ws.Range("N70").Value = -4910.4165

$ws.Range("H73").Value = 7874
$ws.Range("I73").Value = 10877.071
$ws.Range("J73").Value = 4370.4165
$ws.Range("K73").Value = 10877.071
$ws.Range("L73").Value = 4370.4165
$ws.Range("M73").Value = -9941.071
$ws.Range("N73").Value = -6242.4165

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null

$ws.Range("H122").Value = 4764775.5
$ws.Range("I122").Value = 9526323
$ws.Range("J122").Value = 3227.1428
$ws.Range("K122").Value = 28578969
$ws.Range("L122").Value = 9681.428400000001
$ws.Range("M122").Value = -28576519
$ws.Range("N122").Value = -14581.4284

$ws.Range("H132").Value = 3149.5532
$ws.Range("I132").Value = 2211.139
$ws.Range("J132").Value = 6220.727
$ws.Range("K132").Value = 6633.417
$ws.Range("L132").Value = 18662.181
$ws.Range("M132").Value = -4103.417
$ws.Range("N132").Value = -23722.181

$ws.Range("H138").Value = 59228.57
$ws.Range("J138").Value = 59228.57
$ws.Range("L138").Value = 59228.57
$ws.Range("N138").Value = -69508.57000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4925.4
$ws.Range("I40").Value = 4560.2
$ws.Range("J40").Value = 6021
$ws.Range("K40").Value = 4560.2
$ws.Range("L40").Value = 6021
$ws.Range("M40").Value = -4424.2
$ws.Range("N40").Value = -6293

$ws.Range("H132").Value = 7359079
$ws.Range("I132").Value = 3690.5625
$ws.Range("J132").Value = 25012012
$ws.Range("K132").Value = 11071.6875
$ws.Range("L132").Value = 75036036
$ws.Range("M132").Value = -8541.6875
$ws.Range("N132").Value = -75041096

$ws.Range("H136").Value = 10003014
$ws.Range("I136").Value = 10417660
$ws.Range("J136").Value = 51502.5
$ws.Range("K136").Value = 31252980
$ws.Range("L136").Value = 154507.5
$ws.Range("M136").Value = -31250430
$ws.Range("N136").Value = -159607.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 19000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19000
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -19298
